$wb = $excel.ActiveWorkbook

# Helper: force a value to be stored as TEXT (Excel's Value setter otherwise
# auto-detects date-like / numeric-like strings and silently converts them).
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------------
# Sheet 1: "All Published Values" - append a new data row (row 24) for the
# 2026-01-04 publish captured just after the existing 2026-01-03 rows.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

Set-TextValue $ws1.Cells.Item(24, 1) "2026-01-04"
Set-TextValue $ws1.Cells.Item(24, 2) "2026-01-04 00:00:05"
Set-TextValue $ws1.Cells.Item(24, 3) "697.85"
Set-TextValue $ws1.Cells.Item(24, 4) "697.85"
Set-TextValue $ws1.Cells.Item(24, 5) "700.79"
Set-TextValue $ws1.Cells.Item(24, 6) "700.79"
Set-TextValue $ws1.Cells.Item(24, 7) "702.88"
Set-TextValue $ws1.Cells.Item(24, 8) "2026/01/04 00:00:05"
Set-TextValue $ws1.Cells.Item(24, 9) "2026-01-03 16:12:32"
Set-TextValue $ws1.Cells.Item(24, 10) "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"

# Grow the autofilter to cover the new row (A1:J23 -> A1:J24).
$ws1.AutoFilterMode = $false
$ws1.Range("A1:J24").AutoFilter()

# ---------------------------------------------------------------------------
# Sheet 2: "Daily Summary" - fill in today's still-empty summary row (row 6)
# and insert a new "Day First Published" detail row for 2026-01-04.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Row 6 was a blank spacer row directly below the 2026-01-03 averages row;
# populate it in place with the 2026-01-04 averages (single publish so far).
Set-TextValue $ws2.Cells.Item(6, 1) "2026-01-04"
$ws2.Cells.Item(6, 2).Value = 1
$ws2.Cells.Item(6, 3).Value = 702.88
$ws2.Cells.Item(6, 4).Value = 702.88
$ws2.Cells.Item(6, 5).Value = 702.88

# Make room for a new detail row at the bottom of the "Day First Published"
# table by inserting one row just below the old spacer (old row 7), which
# pushes everything from the old row 7 onward down by one row.
$ws2.Rows.Item(7).Insert()

# Populate the newly appended detail row (now row 13) for 2026-01-04.
Set-TextValue $ws2.Cells.Item(13, 1) "2026-01-04"
Set-TextValue $ws2.Cells.Item(13, 2) "2026-01-04 00:00:05"
Set-TextValue $ws2.Cells.Item(13, 3) "702.88"
Set-TextValue $ws2.Cells.Item(13, 4) "2026/01/04 00:00:05"

# The autofilter on this sheet only ever covered the first (averages) table;
# grow it from A1:E5 to A1:E6 to include the new 2026-01-04 averages row.
$ws2.AutoFilterMode = $false
$ws2.Range("A1:E6").AutoFilter()

# ---------------------------------------------------------------------------
# Workbook-level defined names: the hidden _FilterDatabase names mirror each
# sheet's autofilter range and need to be refreshed to match.
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "All Published Values!_FilterDatabase") {
        $n.RefersTo = "='All Published Values'!`$A`$1:`$J`$24"
    }
    elseif ($n.Name -eq "Daily Summary!_FilterDatabase") {
        $n.RefersTo = "='Daily Summary'!`$A`$1:`$E`$6"
    }
}
